$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the results column header from "Z_VNS" to "Z_GVNS" ---
# (Setting the value directly causes the old shared string to be dropped and
#  the new one appended at the end, matching the shared-strings churn seen
#  in the diff.)
$ws.Range("B1").Value = "Z_GVNS"

# --- 2. Clear the old alternating bottom-border style from column B for the
#        "every 10th row" rows (11,21,...,491) so those B-cells fall back to
#        the default (un-bordered) style; column A keeps its border style. ---
for ($r = 11; $r -le 491; $r += 10) {
    $ws.Cells.Item($r, 2).Borders.Item(9).LineStyle = 0
}

# --- 3. Replace the result values for the large (u1000/u2000) instances,
#        rows 402-491, with the new GVNS results, and give that whole block
#        a distinct (black, explicit) font. ---
$newResults = @(514,519,523,531,505,506,505,521,513,505,515,515,527,527,507,510,507,517,512,507,514,518,521,523,515,512,504,520,512,507,509,520,526,524,513,512,510,517,516,507,516,520,525,529,512,518,506,518,515,507,515,521,529,530,516,516,512,524,517,515,518,524,532,535,518,519,512,526,519,515,520,529,533,534,518,524,515,528,520,518,524,531,539,537,523,524,519,529,525,521)

for ($i = 0; $i -lt $newResults.Length; $i++) {
    $ws.Cells.Item(402 + $i, 2).Value = $newResults[$i]
}

$ws.Range("B402:B491").Font.Color = 0

# --- 4. Append 10 brand-new data rows (492-501) that only populate column B
#        (no instance name in column A), using the same new font/style. ---
$extraResults = @(537,540,545,545,527,529,529,535,534,529)

for ($i = 0; $i -lt $extraResults.Length; $i++) {
    $ws.Cells.Item(492 + $i, 2).Value = $extraResults[$i]
}

$ws.Range("B492:B501").Font.Color = 0

# --- 5. Update the selection to match the new, column-B-wide selection seen
#        in the saved file. ---
$ws.Range("B1:B1048576").Select()
